$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.700.22'
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('D3').Value = '1.599.94'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '211.45'
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('E6').Value = '  -0.77%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').Value = '  +0.46%  '
$ws.Range('D9').Value = '0.248'
$ws.Range('E9').Value = '  +1.25%  '
$ws.Range('D10').Value = '19.55'
$ws.Range('E10').Value = '  +0.45%  '
$ws.Range('D11').Value = '0.0841'
$ws.Range('E11').Value = '  +0.70%  '
$ws.Range('D12').Value = '1.824.64'
$ws.Range('E12').Value = '  +0.24%  '
$ws.Range('D13').Value = '1.602.84'
$ws.Range('E13').Value = '  +0.86%  '
$ws.Range('E14').Value = '  +0.69%  '
$ws.Range('D16').Value = '65.34'
$ws.Range('E16').Value = '  +1.27%  '
$ws.Range('D17').Value = '26.676.45'
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('D18').Value = '0.0₃0759'
$ws.Range('E18').Value = '  +3.74%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = '210.11'
$ws.Range('E19').Value = '  +1.07%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '1.00'
$ws.Range('E20').Value = '  +0.12%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').Value = '7.19'
$ws.Range('E21').Value = '  +3.85%  '
$ws.Range('D23').Value = '2.32'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').Value = '8.93'
$ws.Range('E24').Value = '  +0.89%  '
$ws.Range('D25').Value = '143.01'
$ws.Range('E25').Value = '  -1.70%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').Value = '7.12'
$ws.Range('E27').Value = '  -0.11%  '
$ws.Range('E28').Value = '  +0.41%  '
$ws.Range('D29').Value = '15.33'
$ws.Range('E29').Value = '  +0.61%  '
$ws.Range('D30').Value = '0.0517'
$ws.Range('E30').Value = '  +2.70%  '
$ws.Range('E31').Value = '  -0.13%  '
$ws.Range('D32').Value = '3.26'
$ws.Range('E32').Value = '  +1.05%  '
$ws.Range('E33').Value = '  +1.91%  '
$ws.Range('D34').Value = '1.291.56'
$ws.Range('E34').Value = '  +0.74%  '
$ws.Range('D35').Value = '0.619'
$ws.Range('E35').Value = '  -5.23%  '
$ws.Range('D36').Value = '2.47'
$ws.Range('E36').Value = '  +0.81%  '
$ws.Range('E37').Value = '  +0.46%  '
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('D39').Value = '1.07'
$ws.Range('E39').Value = '  +16.08%  '
$ws.Range('E40').Value = '  -2.17%  '
$ws.Range('D41').Value = '5.42'
$ws.Range('E41').Value = '  -0.57%  '
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('E43').Value = '  -0.45%  '
$ws.Range('D44').Value = '63.18'
$ws.Range('E44').Value = '  -1.22%  '
$ws.Range('D45').Value = '1.737.32'
$ws.Range('E45').Value = '  +0.27%  '
$ws.Range('D46').Value = '91.21'
$ws.Range('E46').Value = '  +1.85%  '
$ws.Range('D47').Value = '1.58'
$ws.Range('E47').Value = '  -0.81%  '
$ws.Range('E48').Value = '  -1.18%  '
$ws.Range('E49').Value = '  +0.58%  '
$ws.Range('E50').Value = '  +0.24%  '
$ws.Range('D51').Value = '7.37'
$ws.Range('E51').Value = '  -1.01%  '
